# buff mk9 a little
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 - KAK Value Line AR9
$ws.Range("C11").Value = 9
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 8

# Row 12 - CmmG Mk9 5"
$ws.Range("C12").Value = 7
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 6

# Row 13 - CmmG Mk9 8.5"
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = 3

# Row 14 - CmmG Mk9 9"
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 2

# Row 15 - CmmG Mk9 16" Carbine Length
$ws.Range("C15").Value = -3
$ws.Range("E15").Value = -2
$ws.Range("F15").Value = -2

# Row 16 - CmmG Mk9 16" Rifle Length
$ws.Range("C16").Value = -4
$ws.Range("E16").Value = -3
$ws.Range("F16").Value = -3

# Row 17 - CmmG Mk9 16" Mid Length
$ws.Range("C17").Value = -5
$ws.Range("E17").Value = -4
$ws.Range("F17").Value = -4

# Update selection to match final cursor position in the diff
$ws.Range("C11").Select()
